$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Row 3 on each sheet corresponds to the f0e063cd-... file, now ready for handoff.
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-03-09 09:53:43"

$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-03-09 09:53:46"
